# MAS Report 6A template: clear the example/placeholder figures that were
# left in the "Total value" / "Total number" cells for the digital payment
# token exchange row, so the processor starts from a blank template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 holds "(c) digital payment tokens were exchanged for another
# digital payment token" -- B8 is the total value, C8 the total number.
# Clear their contents but keep the existing number formatting/borders.
$ws.Range("B8:C8").ClearContents()

# Leave the cursor where the author last left it when saving the template.
$ws.Range("C14").Select()
